# Data Driven Automation Using Excel
# Populate the "Actual_Result" column (G) with the Pass/Fail outcome of
# each test row, matching the "Expected_Result" wording style used in F.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "Passed"
$ws.Range("G3").Value = "Passed"
$ws.Range("G4").Value = "Passed"
$ws.Range("G5").Value = "Passed"
$ws.Range("G6").Value = "Failed"

# Give the newly filled cells the plain centered style (no border/fill)
# used elsewhere in the sheet for unshaded, centered content.
$rng = $ws.Range("G2:G6")
$rng.HorizontalAlignment = -4108   # xlCenter
$rng.Borders.LineStyle = -4142     # xlNone
$rng.Interior.Pattern = -4142      # xlNone

$ws.Range("H2:H6").Select()
